$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a trimmed number
# by Excel's type inference (e.g. "1.520" -> 1.52); force Text format first
# so the literal string (matching the source data) is preserved.
$textCells = @('D29', 'D32', 'D33', 'D42', 'D46')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '26.252.80'
$ws.Range('E2').Value = '  -0.07%  '
$ws.Range('D3').Value = '1.663.84'
$ws.Range('E3').Value = '  -0.03%  '
$ws.Range('E4').Value = '  -0.35%  '
$ws.Range('D5').Value = '219.47'
$ws.Range('E5').Value = '  +0.22%  '
$ws.Range('D6').Value = '0.5257'
$ws.Range('E6').Value = '  -0.84%  '
$ws.Range('E7').Value = '  -0.29%  '
$ws.Range('D8').Value = '0.2639'
$ws.Range('E8').Value = '  +0.24%  '
$ws.Range('D9').Value = '0.06325'
$ws.Range('E9').Value = '  -0.54%  '
$ws.Range('D10').Value = '20.75'
$ws.Range('E10').Value = '  +0.91%  '
$ws.Range('D11').Value = '0.07821'
$ws.Range('E11').Value = '  -0.43%  '
$ws.Range('D12').Value = '4.498'
$ws.Range('E12').Value = '  -1.42%  '
$ws.Range('D13').Value = '1.665.72'
$ws.Range('E13').Value = '  +0.02%  '
$ws.Range('D14').Value = '1.891.85'
$ws.Range('E14').Value = '  -0.04%  '
$ws.Range('D15').Value = '0.5561'
$ws.Range('E15').Value = '  +0.57%  '
$ws.Range('D16').Value = '0.0₅8036'
$ws.Range('E16').Value = '  -1.64%  '
$ws.Range('D17').Value = '65.16'
$ws.Range('E17').Value = '  -0.72%  '
$ws.Range('D18').Value = '26.265.80'
$ws.Range('E18').Value = '  -0.13%  '
$ws.Range('E19').Value = '  -0.30%  '
$ws.Range('D20').Value = '4.669'
$ws.Range('E20').Value = '  +0.02%  '
$ws.Range('D21').Value = '197.61'
$ws.Range('E21').Value = '  +2.46%  '
$ws.Range('D22').Value = '10.18'
$ws.Range('E22').Value = '  -0.44%  '
$ws.Range('D23').Value = '5.983'
$ws.Range('E23').Value = '  -1.04%  '
$ws.Range('D24').Value = '1.009'
$ws.Range('E24').Value = '  -0.32%  '
$ws.Range('D25').Value = '145.97'
$ws.Range('E25').Value = '  +0.62%  '
$ws.Range('D26').Value = '0.1211'
$ws.Range('E26').Value = '  -1.05%  '
$ws.Range('D27').Value = '7.183'
$ws.Range('E27').Value = '  -0.43%  '
$ws.Range('D28').Value = '16.07'
$ws.Range('E28').Value = '  -0.51%  '
$ws.Range('D29').Value = '1.520'
$ws.Range('E29').Value = '  +2.09%  '
$ws.Range('D30').Value = '0.05794'
$ws.Range('E30').Value = '  -3.28%  '
$ws.Range('D31').Value = '1.283'
$ws.Range('E31').Value = '  +0.07%  '
$ws.Range('D32').Value = '3.490'
$ws.Range('E32').Value = '  -2.66%  '
$ws.Range('D33').Value = '3.330'
$ws.Range('E33').Value = '  +1.30%  '
$ws.Range('D34').Value = '1.584'
$ws.Range('E34').Value = '  -2.16%  '
$ws.Range('B35').Value = 'MXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D35').Value = '2.815'
$ws.Range('E35').Value = '  -0.48%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').Value = '0.9564'
$ws.Range('E36').Value = '  -0.63%  '
$ws.Range('D37').Value = '2.429'
$ws.Range('E37').Value = '  +0.10%  '
$ws.Range('D38').Value = '0.5775'
$ws.Range('E38').Value = '  -0.45%  '
$ws.Range('D39').Value = '0.01598'
$ws.Range('E39').Value = '  -0.54%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = '5.966'
$ws.Range('E40').Value = '  +0.90%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '1.070.11'
$ws.Range('E41').Value = '  +2.45%  '
$ws.Range('D42').Value = '0.8590'
$ws.Range('E42').Value = '  -0.17%  '
$ws.Range('E43').Value = '  -0.27%  '
$ws.Range('D44').Value = '103.11'
$ws.Range('E44').Value = '  -1.13%  '
$ws.Range('D45').Value = '1.802.65'
$ws.Range('E45').Value = '  -0.14%  '
$ws.Range('D46').Value = '58.50'
$ws.Range('E46').Value = '  +2.00%  '
$ws.Range('B47').Value = 'Frax'
$ws.Range('C47').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D47').Value = '1.011'
$ws.Range('E47').Value = '  -0.63%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').Value = '0.4413'
$ws.Range('E48').Value = '  +0.64%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '7.991'
$ws.Range('E49').Value = '  -0.20%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.05209'
$ws.Range('E50').Value = '  +0.86%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₈101'
$ws.Range('E51').Value = '  -4.95%  '

Write-Output "Applied 114 cell updates"
